# Populate a small daily-driver-report style table on the (single, empty)
# worksheet: a bold/centered/bordered header row followed by one data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel COM constants (this host does not resolve xlCenter/xlTop/xlThin as
# bare names, so the numeric values are used directly):
#   xlCenter = -4108   xlTop = -4160
#   xlContinuous (LineStyle) = 1   xlThin (Weight) = 2
$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1
$xlThin = 2

$headers = @("name", "employee_id", "asset", "arrival", "status", "division", "job_title")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop
    $cell.Borders.LineStyle = $xlContinuous
    $cell.Borders.Weight = $xlThin
}

$values = @("Roger Doddy", "DODROG", "PT-07S", "04:45 AM", "On Time", "TEXDIST", "Select Maintenance Employee")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}
